$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the email for Darshan (row 2): pasne.d -> panse.d
$ws.Range("C2").Value = "panse.d@husky.neu.edu"

# Update the active selection to C2
$ws.Range("C2").Select()
